$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.904.79'
$ws.Range('D3').Value = '1.885.94'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = "'325.43"
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('D7').Value = "'0.4595"
$ws.Range('E7').Value = '  +0.66%  '
$ws.Range('D8').Value = "'0.3894"
$ws.Range('D9').Value = "'0.07849"
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('D10').Value = "'0.9864"
$ws.Range('E10').Value = '  -1.10%  '
$ws.Range('D11').Value = "'21.81"
$ws.Range('E11').Value = '  +0.63%  '
$ws.Range('D12').Value = '1.891.39'
$ws.Range('E12').Value = '  -1.57%  '
$ws.Range('D13').Value = "'7.024"
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('D14').Value = "'5.688"
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D15').Value = "'0.06935"
$ws.Range('E15').Value = '  -0.17%  '
$ws.Range('D16').Value = "'88.25"
$ws.Range('E16').Value = '  +1.53%  '
$ws.Range('D17').Value = "'1.003"
$ws.Range('E17').Value = '  -0.29%  '
$ws.Range('D18').Value = "'0.000009955"
$ws.Range('E18').Value = '  -0.45%  '
$ws.Range('D19').Value = "'17.01"
$ws.Range('E19').Value = '  +0.48%  '
$ws.Range('E20').Value = '  -0.10%  '
$ws.Range('D21').Value = '28.912.45'
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('D22').Value = "'5.280"
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').Value = "'10.98"
$ws.Range('E23').Value = '  +0.35%  '
$ws.Range('D24').Value = '2.111.20'
$ws.Range('E24').Value = '  -1.47%  '
$ws.Range('D25').Value = "'2.085"
$ws.Range('E25').Value = '  +0.60%  '
$ws.Range('D26').Value = "'155.37"
$ws.Range('E26').Value = '  +0.70%  '
$ws.Range('D27').Value = "'19.30"
$ws.Range('E27').Value = '  +0.43%  '
$ws.Range('D28').Value = "'6.006"
$ws.Range('E28').Value = '  +4.63%  '
$ws.Range('D29').Value = "'1.931"
$ws.Range('E29').Value = '  +1.51%  '
$ws.Range('D30').Value = "'117.53"
$ws.Range('E30').Value = '  -0.47%  '
$ws.Range('D31').Value = "'0.09339"
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('D32').Value = "'0.9061"
$ws.Range('E32').Value = '  -0.50%  '
$ws.Range('D33').Value = "'5.282"
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('D34').Value = "'1.328"
$ws.Range('E34').Value = '  +0.14%  '
$ws.Range('D35').Value = "'3.266"
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('D36').Value = "'1.189"
$ws.Range('E36').Value = '  +3.14%  '
$ws.Range('D37').Value = "'0.05763"
$ws.Range('E37').Value = '  +1.12%  '
$ws.Range('D38').Value = "'0.02072"
$ws.Range('E38').Value = '  +1.14%  '
$ws.Range('E39').Value = '  -0.18%  '
$ws.Range('D40').Value = "'7.646"
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').Value = "'0.5671"
$ws.Range('E41').Value = '  +1.51%  '
$ws.Range('D42').Value = "'0.1768"
$ws.Range('E42').Value = '  -0.54%  '
$ws.Range('D43').Value = "'9.698"
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('D44').Value = "'2.253"
$ws.Range('E44').Value = '  +4.99%  '
$ws.Range('D45').Value = "'11.91"
$ws.Range('E45').Value = '  +3.50%  '
$ws.Range('D46').Value = "'0.5352"
$ws.Range('E46').Value = '  +1.38%  '
$ws.Range('E47').Value = '  -1.76%  '
$ws.Range('D48').Value = "'1.850"
$ws.Range('E48').Value = '  +1.92%  '
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').Value = "'2.540"
$ws.Range('E49').Value = '  +3.50%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = "'112.94"
$ws.Range('E50').Value = '  +0.77%  '
$ws.Range('D51').Value = "'1.067"
$ws.Range('E51').Value = '  -4.51%  '
